$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 205
$ws.Range("F3").Value = 114
$ws.Range("F5").Value = 990
$ws.Range("F6").Value = 5490
$ws.Range("F7").Value = 490
$ws.Range("F8").Value = 685
$ws.Range("F10").Value = 822
$ws.Range("F12").Value = 37
$ws.Range("F14").Value = 28
$ws.Range("F17").Value = 1836
$ws.Range("F18").Value = 1469
$ws.Range("F19").Value = 906
$ws.Range("F20").Value = 298
$ws.Range("F22").Value = 331
$ws.Range("F23").Value = 542
$ws.Range("F24").Value = 150
$ws.Range("F25").Value = 1054
$ws.Range("F28").Value = 2872
$ws.Range("F33").Value = 35
$ws.Range("F34").Value = 374
$ws.Range("F36").Value = 42
$ws.Range("F37").Value = 12
$ws.Range("F39").Value = 288
$ws.Range("F40").Value = 710
$ws.Range("F41").Value = 87
$ws.Range("F44").Value = 67
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 34
$ws.Range("F4").Value = 187
$ws.Range("F6").Value = 132
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 237
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 237
$ws.Range("F3").Value = 205
$ws.Range("F4").Value = 114
$ws.Range("F5").Value = 990
$ws.Range("F6").Value = 34
$ws.Range("F7").Value = 5490
$ws.Range("F8").Value = 490
$ws.Range("F9").Value = 685
$ws.Range("F11").Value = 187
$ws.Range("F13").Value = 822
$ws.Range("F15").Value = 132
$ws.Range("F17").Value = 37
$ws.Range("F19").Value = 28
$ws.Range("F23").Value = 1836
$ws.Range("F24").Value = 1469
$ws.Range("F25").Value = 906
$ws.Range("F27").Value = 331
$ws.Range("F29").Value = 542
$ws.Range("F30").Value = 151
$ws.Range("F31").Value = 1054
$ws.Range("F32").Value = 2872
$ws.Range("F37").Value = 35
$ws.Range("F38").Value = 374
$ws.Range("F40").Value = 42
$ws.Range("F41").Value = 12
$ws.Range("F42").Value = 288
$ws.Range("F43").Value = 710
$ws.Range("F44").Value = 87
$ws.Range("F46").Value = 67
